# Apply the edit described by the diff:
# - Sheet "ccterminating": B8 "med" -> "noninversion", B9 "high" -> "inversion",
#                           B11 "med" -> "noninversion", B12 "high" -> "inversion"
# - Select B13 on ccterminating, make it the active/selected sheet (tab)
# - Select C13 on ccplanting (retained, but no longer the active tab)

$wb = $excel.ActiveWorkbook

$wsPlanting = $wb.Worksheets.Item("ccplanting")
$wsTerminating = $wb.Worksheets.Item("ccterminating")

# Update the terminology used on the "ccterminating" sheet
$wsTerminating.Range("B8").Value = "noninversion"
$wsTerminating.Range("B9").Value = "inversion"
$wsTerminating.Range("B11").Value = "noninversion"
$wsTerminating.Range("B12").Value = "inversion"

# Update selections to match the saved view state
$wsPlanting.Range("C13").Select()
$wsTerminating.Range("B13").Select()

# Make "ccterminating" the active sheet/tab
$wsTerminating.Activate()
